# Auto-generated from the cryptos.xlsx OOXML diff.
# Updates Coin/Link/Price/Volume(1h) cells for rows 2-51 of Sheet1.
# Price-column values that look like plain numbers are force-written as
# text (leading apostrophe) so Excel keeps them as strings instead of
# silently coercing them to numeric cells / dropping trailing zeros; the
# style is then reset to Normal so no stray quote-prefix style sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.733.84'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').Value = '2.365.61'
$ws.Range('E3').Value = '  +0.82%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = "'332.64"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.29%  '
$ws.Range('D6').Value = "'101.70"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.58%  '
$ws.Range('E7').Value = '  -0.80%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -0.92%  '
$ws.Range('D10').Value = "'40.08"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.38%  '
$ws.Range('D11').Value = "'0.0923"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.54%  '
$ws.Range('D12').Value = "'8.47"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.06%  '
$ws.Range('E13').Value = '  -3.58%  '
$ws.Range('E14').Value = '  +0.20%  '
$ws.Range('D15').Value = "'16.46"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.14%  '
$ws.Range('D16').Value = '2.726.12'
$ws.Range('E16').Value = '  +0.81%  '
$ws.Range('D17').Value = '2.367.23'
$ws.Range('E17').Value = '  +0.69%  '
$ws.Range('D18').Value = "'8.12"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +12.09%  '
$ws.Range('D19').Value = '42.706.24'
$ws.Range('E19').Value = '  -1.13%  '
$ws.Range('E20').Value = '  -1.49%  '
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').Value = "'76.85"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.03%  '
$ws.Range('B22').Value = 'PancakeSwap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D22').Value = "'3.77"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +9.42%  '
$ws.Range('D23').Value = "'269.87"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.82%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').Value = "'2.33"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -10.00%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').Value = "'10.22"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +12.58%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').Value = "'11.53"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.97%  '
$ws.Range('D28').Value = "'23.23"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.08%  '
$ws.Range('D29').Value = "'2.20"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.14%  '
$ws.Range('D30').Value = "'176.41"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.53%  '
$ws.Range('E31').Value = '  -2.34%  '
$ws.Range('D32').Value = "'0.0901"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.64%  '
$ws.Range('E33').Value = '  -9.43%  '
$ws.Range('D34').Value = "'6.12"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.01%  '
$ws.Range('D35').Value = "'0.132"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  -6.66%  '
$ws.Range('D37').Value = "'2.97"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +10.27%  '
$ws.Range('E38').Value = '  -4.72%  '
$ws.Range('E39').Value = '  +1.07%  '
$ws.Range('D40').Value = "'3.81"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.51%  '
$ws.Range('E41').Value = '  +2.89%  '
$ws.Range('D42').Value = "'0.236"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.60%  '
$ws.Range('D43').Value = "'70.21"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.19%  '
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = "'118.85"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.25%  '
$ws.Range('B46').Value = 'BitcoinSV'
$ws.Range('C46').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D46').Value = "'92.12"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +32.67%  '
$ws.Range('D47').Value = "'11.82"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -7.18%  '
$ws.Range('D48').Value = "'5.52"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.65%  '
$ws.Range('E49').Value = '  -0.98%  '
$ws.Range('E50').Value = '  -2.75%  '
$ws.Range('D51').Value = '1.567.04'
$ws.Range('E51').Value = '  +4.99%  '
